# Upgrade left table: add 2023 column (K) to the Samtredia average monthly
# remuneration table, mirroring the formatting of the existing 2022 column (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New year header in K3 (copy formatting from J3, then set the year value)
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = 2023

# "Gel" row (overall average remuneration) - K4
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 988.9

# "Women" row - K5
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 670.7

# "Men" row - K6
$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("K6").Value = 1299.5999999999999

$excel.CutCopyMode = 0
